$d = $word.ActiveDocument

# Insert a new bold paragraph right after "LAB 2" containing:
#   "Course CODE" <bookmark _GoBack/> ": ETM1142"
$labPara = $d.Paragraphs(1)
$labPara.Range.InsertParagraphAfter()

$codePara = $d.Paragraphs(2)
# Append the first run's text plus a one-character placeholder ("X") so that
# the bookmark insertion point below is not the very last position in the
# paragraph (inserting a collapsed bookmark exactly at paragraph-end is
# unreliable), then strip the placeholder back out afterwards.
$codePara.Range.InsertAfter("Course CODEX")

$bmPos = $codePara.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($codePara.Range.End - 2, $codePara.Range.End - 1)
$placeholder.Text = ": ETM1142"
